# Add a new "Cost (£)/pu (estimate)" column (F) to the TV consumption
# statistics sheet, with unit-price figures for each TV type.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header in F1, matching the style of the other header cells (B1:E1, G1)
$ws.Range("F1").Value = "Cost (£)/pu (estimate)"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row needs to grow to fit the wrapped new heading text
$ws.Rows.Item(1).RowHeight = 42

# Per-row unit price estimates
$prices = @{
    2  = 280
    3  = 420
    4  = 900
    5  = 250
    6  = 320
    7  = 380
    8  = 700
    9  = 1200
    10 = 2000
    11 = 2500
    12 = 300
    13 = 300
    14 = 350
    15 = 450
    16 = 1000
    17 = 2000
    18 = 2500
    19 = 50
}

foreach ($row in $prices.Keys) {
    $ws.Cells.Item($row, 6).Value = $prices[$row]
}

# Move the active selection down to reflect the newly-populated column
$ws.Range("F20").Select()
